{"js": "// Update the date line and the 25 division equations in the practice table.\n// The table has 20 rows x 5 columns, but only every 4th row (0, 4, 8, 12, 16)\n// actually holds equation text -- the rows in between are blank answer rows.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Update the date/weekday heading (first paragraph of the document).\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\nif (dateParagraph.text.indexOf(\"2023-07-22 Saturday\") !== -1) {\n  dateParagraph.getRange().insertText(\"2023-07-23 Sunday\", \"Replace\");\n}\n\n// 2) Update the table of two-digit division problems.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// old -> new text for each of the 25 equation cells, given in row-major\n// order (5 rows of 5 columns each) matching the table's content rows.\nconst replacements = [\n  [\"44\u00f73=\", \"77\u00f72=\"],\n  [\"31\u00f76=\", \"22\u00f74=\"],\n  [\"42\u00f75=\", \"77\u00f77=\"],\n  [\"13\u00f79=\", \"53\u00f73=\"],\n  [\"78\u00f79=\", \"90\u00f76=\"],\n\n  [\"39\u00f74=\", \"30\u00f73=\"],\n  [\"85\u00f78=\", \"81\u00f79=\"],\n  [\"20\u00f79=\", \"71\u00f78=\"],\n  [\"38\u00f78=\", \"16\u00f77=\"],\n  [\"61\u00f77=\", \"62\u00f79=\"],\n\n  [\"89\u00f77=\", \"27\u00f75=\"],\n  [\"96\u00f74=\", \"41\u00f73=\"],\n  [\"74\u00f77=\", \"21\u00f76=\"],\n  [\"78\u00f73=\", \"21\u00f72=\"],\n  [\"32\u00f77=\", \"66\u00f72=\"],\n\n  [\"52\u00f75=\", \"84\u00f73=\"],\n  [\"54\u00f75=\", \"80\u00f78=\"],\n  [\"44\u00f73=\", \"11\u00f72=\"],\n  [\"82\u00f76=\", \"35\u00f75=\"],\n  [\"16\u00f72=\", \"57\u00f73=\"],\n\n  [\"15\u00f79=\", \"53\u00f75=\"],\n  [\"55\u00f76=\", \"63\u00f73=\"],\n  [\"78\u00f78=\", \"74\u00f76=\"],\n  [\"37\u00f79=\", \"91\u00f76=\"],\n  [\"77\u00f75=\", \"12\u00f76=\"],\n];\n\n// The 5 content rows of the table (0-indexed) that hold equations.\nconst contentRows = [0, 4, 8, 12, 16];\nconst columnsPerRow = 5;\n\nfor (let r = 0; r < contentRows.length; r++) {\n  const tableRow = contentRows[r];\n  for (let c = 0; c < columnsPerRow; c++) {\n    const idx = r * columnsPerRow + c;\n    const [oldText, newText] = replacements[idx];\n    const cell = table.getCell(tableRow, c);\n    cell.load(\"value\");\n    await context.sync();\n    if (cell.value.indexOf(oldText) !== -1) {\n      cell.body.getRange().insertText(newText, \"Replace\");\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division equations in the practice table.\n# The table has 20 rows x 5 columns, but only every 4th row (COM rows\n# 1, 5, 9, 13, 17) actually holds equation text -- the rows in between are\n# blank answer rows for students to fill in.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/weekday heading (first paragraph of the document).\n$dateRange = $d.Paragraphs(1).Range\nif ($dateRange.Text -like \"*2023-07-22 Saturday*\") {\n    $dateRange.Text = \"2023-07-23 Sunday\"\n}\n\n# 2) Update the table of two-digit division problems.\n$table = $d.Tables(1)\n\n# COM table rows (1-indexed) that hold equation text.\n$contentRows = @(1, 5, 9, 13, 17)\n\n# old -> new text for each of the 25 equation cells, given in row-major\n# order (5 rows of 5 columns each) matching the table's content rows.\n$replacements = @(\n    @(\"44\u00f73=\", \"77\u00f72=\"),\n    @(\"31\u00f76=\", \"22\u00f74=\"),\n    @(\"42\u00f75=\", \"77\u00f77=\"),\n    @(\"13\u00f79=\", \"53\u00f73=\"),\n    @(\"78\u00f79=\", \"90\u00f76=\"),\n\n    @(\"39\u00f74=\", \"30\u00f73=\"),\n    @(\"85\u00f78=\", \"81\u00f79=\"),\n    @(\"20\u00f79=\", \"71\u00f78=\"),\n    @(\"38\u00f78=\", \"16\u00f77=\"),\n    @(\"61\u00f77=\", \"62\u00f79=\"),\n\n    @(\"89\u00f77=\", \"27\u00f75=\"),\n    @(\"96\u00f74=\", \"41\u00f73=\"),\n    @(\"74\u00f77=\", \"21\u00f76=\"),\n    @(\"78\u00f73=\", \"21\u00f72=\"),\n    @(\"32\u00f77=\", \"66\u00f72=\"),\n\n    @(\"52\u00f75=\", \"84\u00f73=\"),\n    @(\"54\u00f75=\", \"80\u00f78=\"),\n    @(\"44\u00f73=\", \"11\u00f72=\"),\n    @(\"82\u00f76=\", \"35\u00f75=\"),\n    @(\"16\u00f72=\", \"57\u00f73=\"),\n\n    @(\"15\u00f79=\", \"53\u00f75=\"),\n    @(\"55\u00f76=\", \"63\u00f73=\"),\n    @(\"78\u00f78=\", \"74\u00f76=\"),\n    @(\"37\u00f79=\", \"91\u00f76=\"),\n    @(\"77\u00f75=\", \"12\u00f76=\")\n)\n\nfor ($r = 0; $r -lt $contentRows.Length; $r++) {\n    $tableRow = $contentRows[$r]\n    for ($c = 1; $c -le 5; $c++) {\n        $idx = $r * 5 + ($c - 1)\n        $pair = $replacements[$idx]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n        $cell = $table.Cell($tableRow, $c)\n        if ($cell.Range.Text -like \"*$oldText*\") {\n            $cell.Range.Text = $newText\n        }\n    }\n}\n"}
